$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brix_Gel_Stab")
$ws.Activate()

# Fix the "Handlauf" options text: swap order of the two choices
$ws.Range("D5").Value = "Freilauf Handlauf:1.25, Integrierter Handlauf:1"

# Insert two new rows before the current row 9 ("Montage (€/m)") to make
# room for the new "Zier-Element" parameters, pushing the old rows 9-11 down
# to rows 11-13.
$ws.Rows("9:10").Insert()

# New row 9: Zier-Element selector (Auswahl)
$ws.Range("A9").Value = "Auswahl"
$ws.Range("B9").Value = "Zier-Element"
$ws.Range("C9").Value = "P_Zier"
$ws.Range("D9").Value = "NEIN:0, Ja:50"

# New row 10: Zier-Element Anzahl (Zahl)
$ws.Range("A10").Value = "Zahl"
$ws.Range("B10").Value = "Zier-Element Anzahl"
$ws.Range("C10").Value = "Zier_Stk"

# Update the final price formula (now on row 13, column E) to include the
# new Zier-Element contribution.
$ws.Range("E13").Value = "((P_Modell * L * F_Faktor * P_Handlauf) + ((math.ceil(L/1.3)+1) * P_Steher * F_Faktor) + (Ecken * 95) + (L * P_Arbeit) + (L * F_Schräg) + (P_Zier * Zier_Stk)) * ( 1 - (p_rabatt / 100))"

$ws.Range("E13").Select()
